$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data from source feed
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.916.54'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.78%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.721.92'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '609.56'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '183.32'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.93%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.718.32'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.633'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.63%  '
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.723'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.87%  '
$ws.Range('E11').Value = '  -8.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '57.47'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.99%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000296'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.73'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.308.59'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.721.36'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.51'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.46%  '
$ws.Range('E18').Value = '  -1.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.99'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.08%  '
$ws.Range('E20').Value = '  -6.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '68.760.79'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '415.41'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.73'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '89.31'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.45%  '
$ws.Range('E25').Value = '  -7.75%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.78'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -8.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.57%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.90'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.08'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('E30').Value = '  -8.37%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.15'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.37'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -15.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.56'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.84%  '
$ws.Range('E34').Value = '  -5.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '44.15'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '65.28'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.64%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '607.84'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0893'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -11.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.409'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.35%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('E42').Value = '  -5.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.08'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0444'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.68'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.32%  '
$ws.Range('E46').Value = '  -11.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.26'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -9.73%  '
$ws.Range('E48').Value = '  -6.32%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.795.31'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.15%  '
$ws.Range('B50').Value = 'WEMIXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.72'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.73%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.07'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.92%  '
